$d = $word.ActiveDocument

function Find-Range([string]$text) {
    # Locate $text in the document and return the Range collapsed onto
    # the match (Start/End reflect the match position), like real Word.
    $r = $d.Content
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $text"
    }
    return $r
}

function Split-Into-Runs([string]$newText, [string[]]$parts) {
    # $newText must already be present verbatim in the document (usually
    # just placed there via a Find/Replace). Re-locate it and force run
    # boundaries at each internal join between consecutive $parts,
    # without altering the visible text/formatting: toggle a
    # character property on from the cut point to the end of the match,
    # then back off. Because both the "on" and "off" operations share
    # the exact same start boundary, this yields one extra run per cut
    # (not one per character).
    $rng = Find-Range $newText
    $base = $rng.Start
    $end = $rng.End

    $off = 0
    $cuts = @()
    foreach ($p in $parts) {
        $off += $p.Length
        $cuts += $off
    }
    # Walk cut points right-to-left so earlier (still-to-be-split)
    # prefix text is always covered by a single Find per call.
    for ($i = $cuts.Length - 2; $i -ge 0; $i--) {
        $pos = $base + $cuts[$i]
        $rr = $d.Range($pos, $end)
        $rr.Font.Bold = 1
        $rr.Font.Bold = 0
    }
}

# ---------------------------------------------------------------------
# Change 1 (home page / header paragraph): "Used grid to keep header
# columns as window is scaled." -> split into 5 runs, "grid" becomes
# "flexbox" and "columns" becomes "formatted as".
# ---------------------------------------------------------------------
$old1 = " Used grid to keep header columns as window is scaled."
$new1 = " Used flexbox to keep header formatted as window is scaled."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null
Split-Into-Runs $new1 @(" Used ", "flexbox", " to keep header ", "formatted as", " window is scaled.")

# ---------------------------------------------------------------------
# Change 2 (logo/banner paragraph): "...if I try to make it max-hight
# 80vh. " -> "...if I try to adjust its max-height. ", split into 7
# runs.
# ---------------------------------------------------------------------
$old2 = "Due to not using full version of logo as icon, wanted to add larger one top of homepage which would also add to the user flow. Image gets stretched (at least on my screen) if I try to make it max-hight 80vh. "
$new2 = "Due to not using full version of logo as icon, wanted to add larger one top of homepage which would also add to the user flow. Image gets stretched (at least on my screen) if I try to adjust its max-height. "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null
Split-Into-Runs $new2 @(
    "Due to not using full version of logo as icon, wanted to add larger one top of homepage which would also add to the user flow. Image gets stretched (at least on my screen) if I try to ",
    "adjust",
    " it",
    "s",
    " max-h",
    "e",
    "ight. "
)

# ---------------------------------------------------------------------
# Change 3 (viewport paragraph): append a new trailing run about
# px vs vw.
# ---------------------------------------------------------------------
$p3 = Find-Range "Had to adjust some of my widths to vw as once viewport started working, my scaling was no longer responsive."
$p3.Collapse(0) | Out-Null
$p3.InsertAfter(" A little confused as to if I should be using a fixed unit such as px or responsive such as vw. Research suggests vw, but requires further research to gain more understanding.") | Out-Null

# ---------------------------------------------------------------------
# Change 4 (nav bar): insert a new paragraph "Used flexbox for nav
# bar..." right before "Added hover and active colours..." paragraph.
# ---------------------------------------------------------------------
$hover = Find-Range "Added hover and active colours to the anchors in the nav bar."
$hoverPara = $hover.Paragraphs.Item(1)
$hoverPara.Range.InsertBefore("Used flexbox for nav bar. Allowing hyperlinks to be side by side on desktop, then on top of each other for mobile.`r")

# ---------------------------------------------------------------------
# Change 5: after the final "Made nav bar sticky..." paragraph, add two
# blank paragraphs, a TODO-style paragraph, then one more blank
# paragraph.
# ---------------------------------------------------------------------
$sticky = Find-Range "Made nav bar sticky, so it gets moved up until just under the header bar and becomes fixed."
$stickyPara = $sticky.Paragraphs.Item(1)
$stickyPara.Range.InsertParagraphAfter()
$afterSticky = $d.Paragraphs.Item($stickyPara.Index + 1)
$afterSticky.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($stickyPara.Index + 2)
$p2.Range.InsertParagraphAfter()
$p3todo = $d.Paragraphs.Item($stickyPara.Index + 3)
$p3todo.Range.InsertAfter("ADJUST NAV BAR STIKY POSITION FOR MOBILE ")
$p3todo.Range.InsertParagraphAfter()

Write-Output "done"
